# Form submission flow case
# Adds a new "1099MISCdata" worksheet with payer business-name test data,
# tweaks a couple of stored selections, and normalizes the EIN sheet's
# numeric value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "1099MISCdata" worksheet as the last tab (after
#    "AddRecipientEin"), matching sheetId=17 / rId16 placement.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "1099MISCdata"

# Column widths (character units -> stored width), closest achievable
# quantization to 27.88671875 / 21.6640625.
$newSheet.Columns.Item(1).ColumnWidth = 27.0
$newSheet.Columns.Item(2).ColumnWidth = 20.833333333333332

# Header / label cells
$newSheet.Range("A1").Value = "AL"
$newSheet.Range("A2").Value = "R123456789"
$newSheet.Range("A3").Value = "BusinessName"
$newSheet.Range("A4").Value = "Payer Business name"

# The test-data values in column B went through a few intermediate
# entries while the tester tried different generated IDs before the
# final values were entered.
$newSheet.Range("B3").Value = "Test0202202120528"
$newSheet.Range("B3").Value = ""
$newSheet.Range("B4").Value = "Test0202202134948"
$newSheet.Range("B4").Value = "Test0202202150222"
$newSheet.Range("B3").Value = "Test0202202155851"
$newSheet.Range("B3").Value = "Test0202202160448"

$newSheet.Range("B4").Value = "Test0203202195509"
$newSheet.Range("B3").Value = "Test02032021100108"

# ---------------------------------------------------------------------
# 2. EIN sheet: re-enter the payer id so it is stored as a plain
#    integer value.
# ---------------------------------------------------------------------
$einSheet = $wb.Worksheets.Item("EIN")
$einSheet.Range("A1").Value = 546199841

# ---------------------------------------------------------------------
# 3. AddPayerNonUsEin: selection moved from F3 to G3.
# ---------------------------------------------------------------------
$nonUsEinSheet = $wb.Worksheets.Item("AddPayerNonUsEin")
$nonUsEinSheet.Range("G3").Select()

# ---------------------------------------------------------------------
# 4. AddRecipientEin: selection moved from C7 to H13 (it also stops
#    being the active/saved tab once the new sheet becomes active).
# ---------------------------------------------------------------------
$addRecipientEinSheet = $wb.Worksheets.Item("AddRecipientEin")
$addRecipientEinSheet.Range("H13").Select()

# ---------------------------------------------------------------------
# 5. Leave the newly added sheet active/selected at B4, as the last
#    worked-on sheet.
# ---------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("B4").Select()
